$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new columns before the old "L" column (right after I), which pushes
#    ListObject2 / NamedRange2 / ListObject3 / NamedRange3 content 3 columns to the right.
$ws.Columns("J:L").Insert()

# 2. Fix up the table ranges that the column insert does not auto-adjust in this runtime.
$lo2 = $ws.ListObjects.Item("ListObject2")
$lo2.Resize($ws.Range("O4:Q5"))
$lo3 = $ws.ListObjects.Item("ListObject3")
$lo3.Resize($ws.Range("Y4:Y5"))

# 3. Fix up the defined names similarly.
$wb.Names.Item("NamedRange1").RefersTo = "=TableToDicts!`$G`$4:`$L`$6"
$wb.Names.Item("NamedRange2").RefersTo = "=TableToDicts!`$T`$4:`$V`$4"
$wb.Names.Item("NamedRange3").RefersTo = "=TableToDicts!`$AB`$4"

# 4. Rearrange the content inside the (now widened) NamedRange1 block G4:L6 so it
#    demonstrates duplicate / blank column headers (bugfix: allow duplicate column names).
#    Move old column H (b/2/5) into I, blank out H, and populate new K (C/30/60) and L (c/3/6).
$ws.Range("I4").Value = "b"
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 5

$ws.Range("H4").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("H6").ClearContents()

$ws.Range("K4").Value = "C"
$ws.Range("K5").Value = 30
$ws.Range("K6").Value = 60

$ws.Range("L4").Value = "c"
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 6

# 5. Note above the duplicate-name columns.
$ws.Range("K3").Value = "Case insensitive and uses last instance of c"

# 6. Shade the whole NamedRange1 block (G4:L6) with a light gray fill, like a table header.
#    (Cell-by-cell on purpose: applying ThemeColor to a multi-cell Range in a single call
#    does not reliably commit the themed fill to every cell in this runtime.)
foreach ($colLetter in @("G", "H", "I", "J", "K", "L")) {
    foreach ($rowNum in 4..6) {
        $cell = $ws.Range("$colLetter$rowNum")
        $cell.Interior.ThemeColor = 2
        $cell.Interior.TintAndShade = -0.14999847407452621
    }
}

# 7. Selection / view tidy-up to match the target sheet view.
$ws.Range("H9").Select()
